$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8811536431312561
$ws.Range("B1").Value = 2.703952074050903
$ws.Range("C1").Value = 3.384826898574829
$ws.Range("D1").Value = 1.913663983345032
$ws.Range("E1").Value = 1.468409180641174
